$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-10 Tuesday" "2026-02-11 Wednesday"

Replace-Text "92÷9=" "82÷9="
Replace-Text "48÷5=" "49÷4="
Replace-Text "95÷9=" "49÷6="
Replace-Text "23÷5=" "37÷4="
Replace-Text "75÷8=" "99÷2="
Replace-Text "19÷6=" "87÷3="
Replace-Text "73÷5=" "38÷7="
Replace-Text "14÷5=" "87÷9="
Replace-Text "72÷3=" "69÷9="
Replace-Text "63÷3=" "92÷6="
Replace-Text "71÷2=" "49÷2="
Replace-Text "78÷6=" "36÷9="
Replace-Text "22÷9=" "86÷3="
Replace-Text "45÷7=" "15÷2="
Replace-Text "16÷2=" "38÷6="
Replace-Text "48÷9=" "12÷2="
Replace-Text "41÷9=" "69÷2="
Replace-Text "38÷9=" "26÷5="
Replace-Text "95÷8=" "68÷6="
Replace-Text "53÷6=" "27÷3="
Replace-Text "26÷4=" "72÷3="
Replace-Text "17÷6=" "49÷3="
Replace-Text "19÷4=" "83÷2="
Replace-Text "30÷5=" "68÷9="
Replace-Text "22÷2=" "64÷4="

"Done"
